$wb = $excel.ActiveWorkbook

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6629.4614
$ws.Range("I32").Value = 6629.4614
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 6629.4614
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -6342.4614
$ws.Range("N32").ClearContents()
# Row 61
$ws.Range("H61").Value = 3159.6
$ws.Range("I61").Value = 2955.111
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2955.111
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2743.111
$ws.Range("N61").Value = -5424
# Row 136
$ws.Range("H136").Value = 3159.6
$ws.Range("I136").Value = 2955.111
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 8865.332999999999
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -6315.332999999999
$ws.Range("N136").Value = -20100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3542
$ws.Range("I20").Value = 3278.3333
$ws.Range("J20").Value = 4333
$ws.Range("K20").Value = 3278.3333
$ws.Range("L20").Value = 4333
$ws.Range("M20").Value = -3031.3333
$ws.Range("N20").Value = -4827
# Row 22
$ws.Range("H22").Value = 179
$ws.Range("I22").Value = 134.8
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 134.8
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = 38.19999999999999
$ws.Range("N22").Value = -746
# Row 94
$ws.Range("H94").Value = 374.5
$ws.Range("I94").Value = 374.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 374.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 76.5
$ws.Range("N94").ClearContents()
# Row 134
$ws.Range("H134").Value = 2117.647
$ws.Range("I134").Value = 1937.5
$ws.Range("K134").Value = 5812.5
$ws.Range("M134").Value = -3277.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7706.6553
$ws.Range("I31").Value = 4191.4287
$ws.Range("J31").Value = 8825.137000000001
$ws.Range("K31").Value = 4191.4287
$ws.Range("L31").Value = 8825.137000000001
$ws.Range("M31").Value = -3896.4287
$ws.Range("N31").Value = -9415.137000000001
# Row 34
$ws.Range("H34").Value = 7706.6553
$ws.Range("I34").Value = 4191.4287
$ws.Range("J34").Value = 8825.137000000001
$ws.Range("K34").Value = 4191.4287
$ws.Range("L34").Value = 8825.137000000001
$ws.Range("M34").Value = -3989.4287
$ws.Range("N34").Value = -9229.137000000001
# Row 134
$ws.Range("H134").Value = 2075.5
$ws.Range("I134").Value = 2075.5
$ws.Range("K134").Value = 6226.5
$ws.Range("M134").Value = -3691.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 407.8095
$ws.Range("J17").Value = 940.75
$ws.Range("L17").Value = 2822.25
$ws.Range("N17").Value = -3160.25
# Row 68
$ws.Range("H68").Value = 626
$ws.Range("J68").Value = 839
$ws.Range("L68").Value = 2517
$ws.Range("N68").Value = -4139
# Row 71
$ws.Range("H71").Value = 626
$ws.Range("J71").Value = 839
$ws.Range("L71").Value = 7551
$ws.Range("N71").Value = -15663
# Row 109
$ws.Range("H109").Value = 2700
$ws.Range("I109").Value = 1566.6666
$ws.Range("K109").Value = 4699.9998
$ws.Range("M109").Value = -3659.9998
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
# Row 140
$ws.Range("H140").Value = 2102.1538
$ws.Range("I140").Value = 1529.8182
$ws.Range("K140").Value = 4589.4546
$ws.Range("M140").Value = 590.5454
# Row 141
$ws.Range("H141").Value = 2496.25
$ws.Range("I141").Value = 2496.25
$ws.Range("K141").Value = 7488.75
$ws.Range("M141").Value = -2308.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 20000
$ws.Range("I57").Value = 20000
$ws.Range("K57").Value = 20000
$ws.Range("M57").Value = -19180
# Row 80
$ws.Range("H80").Value = 1003
$ws.Range("I80").Value = 1003
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1003
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -5
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 1003
$ws.Range("I83").Value = 1003
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 5015
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -23
$ws.Range("N83").ClearContents()
# Row 132
$ws.Range("H132").Value = 29743.309
$ws.Range("I132").Value = 38110.9
$ws.Range("J132").Value = 5477.3
$ws.Range("K132").Value = 114332.7
$ws.Range("L132").Value = 16431.9
$ws.Range("M132").Value = -111802.7
$ws.Range("N132").Value = -21491.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2993.75
$ws.Range("I16").Value = 658.6667
$ws.Range("K16").Value = 658.6667
$ws.Range("M16").Value = -488.6667
# Row 68
$ws.Range("H68").Value = 5677.6
$ws.Range("I68").Value = 4597
$ws.Range("K68").Value = 4597
$ws.Range("M68").Value = -3848
# Row 71
$ws.Range("H71").Value = 5677.6
$ws.Range("I71").Value = 4597
$ws.Range("K71").Value = 22985
$ws.Range("M71").Value = -19241
# Row 82
$ws.Range("H82").Value = 6648.625
$ws.Range("I82").Value = 5297.25
$ws.Range("K82").Value = 5297.25
$ws.Range("M82").Value = -4936.25
# Row 85
$ws.Range("H85").Value = 6648.625
$ws.Range("I85").Value = 5297.25
$ws.Range("K85").Value = 5297.25
$ws.Range("M85").Value = -4049.25
# Row 132
$ws.Range("H132").Value = 5337
$ws.Range("I132").Value = 4559
$ws.Range("J132").Value = 10005
$ws.Range("K132").Value = 13677
$ws.Range("L132").Value = 30015
$ws.Range("M132").Value = -11147
$ws.Range("N132").Value = -35075
# Row 136
$ws.Range("H136").Value = 4161
$ws.Range("I136").Value = 2518.3333
$ws.Range("K136").Value = 7554.999899999999
$ws.Range("M136").Value = -5004.999899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 17
$ws.Range("H17").Value = 2502
$ws.Range("I17").Value = 2004
$ws.Range("J17").Value = 3000
$ws.Range("K17").Value = 2004
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -1832
$ws.Range("N17").Value = -3344
# Row 81
$ws.Range("H81").Value = 1330.6666
$ws.Range("I81").Value = 1495
$ws.Range("J81").Value = 1002
$ws.Range("K81").Value = 2990
$ws.Range("L81").Value = 2004
$ws.Range("M81").Value = -1929
$ws.Range("N81").Value = -4126
# Row 84
$ws.Range("H84").Value = 1330.6666
$ws.Range("I84").Value = 1495
$ws.Range("J84").Value = 1002
$ws.Range("K84").Value = 14950
$ws.Range("L84").Value = 10020
$ws.Range("M84").Value = -9646
$ws.Range("N84").Value = -20628
# Row 113
$ws.Range("H113").Value = 1290.4
$ws.Range("I113").Value = 1290.4
$ws.Range("K113").Value = 3871.2
$ws.Range("M113").Value = -1701.2
